$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $old"
    }
    return $result
}

# AstEvaluator.caseCall line number change 189 -> 186
Replace-Text "at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)" "at org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:186)"

# AstEvaluator.eval line number change 112 -> 109
Replace-Text "at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)" "at org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)"

# M2DocEvaluator.caseUserDoc line number change 1252 -> 1270
Replace-Text "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1252)" "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1270)"

# doSwitch(1216) + caseBlock(1425) -> doSwitch(1234) + caseBlock(1459)
$old1 = "doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)"
$new1 = "doSwitch(M2DocEvaluator.java:1234)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1459)"
Replace-Text $old1 $new1

# doSwitch(1216) + caseDocumentTemplate(287) -> doSwitch(1234) + caseDocumentTemplate(296)
$old2 = "doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)"
$new2 = "doSwitch(M2DocEvaluator.java:1234)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)"
Replace-Text $old2 $new2

# doSwitch(1216) + generate(276) -> doSwitch(1234) + generate(281)
$old3 = "doSwitch(M2DocEvaluator.java:1216)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:276)"
$new3 = "doSwitch(M2DocEvaluator.java:1234)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:281)"
Replace-Text $old3 $new3

# M2DocUtils.generate line number change 694 -> 696
Replace-Text "at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:694)" "at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:696)"

# AbstractTemplatesTestSuite.prepareoutputAndGenerate line number change 480 -> 463
Replace-Text "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)" "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:463)"

# AbstractTemplatesTestSuite.generation line number change 389 -> 373
Replace-Text "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)" "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:373)"

# GeneratedMethodAccessor75 line replaced with two NativeMethodAccessorImpl lines
$old4 = "at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"
$new4 = "at sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)"
Replace-Text $old4 $new4

# Remove duplicated 8-line block of ParentRunner/Suite trace entries
$delOld = "`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)"
Replace-Text $delOld ""
